# formula: support for more math/trig formulas
# Adds columns U:AG on the "Math and Trig" sheet for:
# EVEN, EXP, FACT, FACTDOUBLE, FLOOR.MATH, FLOOR.PRECISE, GCD, INT,
# ISO.CEILING, LCM, LN, LOG, LOG10
# mirroring the layout/style already used for the other math functions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- Row 2: bold, centered headers (matches style used by B2:T2) ----
$headers = @{
    "U2"  = "EVEN"
    "V2"  = "EXP"
    "W2"  = "FACT"
    "X2"  = "FACTDOUBLE"
    "Y2"  = "FLOOR.MATH"
    "Z2"  = "FLOOR.PRECISE"
    "AA2" = "GCD"
    "AB2" = "INT"
    "AC2" = "ISO.CEILING"
    "AD2" = "LCM"
    "AE2" = "LN"
    "AF2" = "LOG"
    "AG2" = "LOG10"
}
foreach ($ref in $headers.Keys) {
    $cell = $ws.Range($ref)
    $cell.Value = $headers[$ref]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

# ---- Formula cells, row by row ----
$formulas = @{
    # Row 3
    "U3"  = "=EVEN()"
    "V3"  = "=EXP(U3)"
    "W3"  = "=FACT()"
    "X3"  = "=FACTDOUBLE()"
    "Y3"  = "=FLOOR.MATH()"
    "Z3"  = "=FLOOR.PRECISE()"
    "AA3" = "=GCD()"
    "AB3" = "=INT()"
    "AC3" = "=ISO.CEILING()"
    "AD3" = "=LCM()"
    "AE3" = "=LN()"
    "AF3" = "=LOG()"
    "AG3" = "=LOG10()"

    # Row 4
    "U4"  = "=EVEN(1.6)"
    "V4"  = "=EXP()"
    "W4"  = "=FACT(-1)"
    "X4"  = "=FACTDOUBLE(0)"
    "Y4"  = "=FLOOR.MATH(3.7,2)"
    "Z4"  = "=FLOOR.PRECISE(-3.2,-1)"
    "AA4" = "=GCD(Z4:Z7)"
    "AB4" = "=INT(0)"
    "AC4" = "=ISO.CEILING(0)"
    "AD4" = "=LCM(5,2)"
    "AE4" = "=LEFT(LN(86),5)"
    "AF4" = "=LOG(10)"
    "AG4" = "=LOG10(100)"

    # Row 5
    "U5"  = "=EVEN(3)"
    "V5"  = "=LEFT(EXP(U5),5)"
    "W5"  = "=FACT(0)"
    "X5"  = "=FACTDOUBLE(-1)"
    "Y5"  = "=FLOOR.MATH(-2.5,-2)"
    "Z5"  = "=FLOOR.PRECISE(3.2,1)"
    "AA5" = "=GCD(Y6:Y8)"
    "AB5" = "=INT(1)"
    "AC5" = "=ISO.CEILING(4.3)"
    "AD5" = "=LCM(24,36)"
    "AE5" = "=LN(EXP(3.15))"
    "AF5" = "=LOG(8,2)"
    "AG5" = "=LOG10(0)"

    # Row 6
    "U6"  = "=EVEN(2)"
    "V6"  = "=LEFT(EXP(U7),5)"
    "W6"  = "=FACT(1)"
    "X6"  = "=FACTDOUBLE(6)"
    "Y6"  = "=FLOOR.MATH(2.5,-2)"
    "Z6"  = "=FLOOR.PRECISE(-3.2,1)"
    "AA6" = "=GCD(5,2)"
    "AB6" = "=INT(8.9)"
    "AC6" = "=ISO.CEILING(-4.3)"
    "AD6" = "=LCM(1,2,3,5)"
    "AE6" = "=LN(EXP(3))"
    "AF6" = "=LOG(1024,2)"
    "AG6" = "=LOG10(1000)"

    # Row 7
    "U7"  = "=EVEN(-1)"
    "W7"  = "=FACT(5)"
    "X7"  = "=FACTDOUBLE(7)"
    "Y7"  = "=FLOOR.MATH(1.58,0.1)"
    "Z7"  = "=FLOOR.PRECISE(3.2,-1)"
    "AA7" = "=GCD(24,36)"
    "AB7" = "=INT(-8.9)"
    "AC7" = "=ISO.CEILING(4.3,2)"
    "AD7" = "=LCM(0)"
    "AF7" = "=LOG(0)"

    # Row 8
    "U8"  = "=EVEN(U6:U7)"
    "Y8"  = "=FLOOR.MATH(0.234,0.01)"
    "Z8"  = "=FLOOR.PRECISE(3.2)"
    "AA8" = "=GCD(7,1)"
    "AB8" = "=INT(0.5)"
    "AC8" = "=ISO.CEILING(4.3,-2)"
    "AD8" = "=LCM(-1)"
    "AF8" = "=LOG(10,0)"

    # Row 9
    "Y9"  = "=FLOOR.MATH(1.5)"
    "AA9" = "=GCD(5,0)"
    "AB9" = "=INT(-0.5)"
    "AC9" = "=ISO.CEILING(-4.3,2)"
    "AD9" = "=LCM(5,2,-1)"
    "AF9" = "=LOG(1,0)"

    # Row 10
    "Y10"  = "=FLOOR.MATH(-1.5)"
    "AA10" = "=GCD(5)"
    "AB10" = "=INT(0)"
    "AC10" = "=ISO.CEILING(-4.3,-2)"
    "AD10" = "=LCM(5,23,0)"

    # Row 11
    "AA11" = '=GCD(5,8,"a")'
    "AD11" = "=LCM(0,0)"

    # Row 12
    "AA12" = "=GCD(2,0)"

    # Row 13
    "AA13" = "=GCD(0,5)"

    # Row 14
    "AA14" = "=GCD(0,0,6)"
}
foreach ($ref in $formulas.Keys) {
    $ws.Range($ref).Formula = $formulas[$ref]
}

# ---- Window/selection state (matches the commit's saved view) ----
$ws.Range("AG7").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 14
